# Add two new Mac-Addresses (10 new device rows) to the reg_center/machine/device
# master data sheet, mirroring the existing row pattern but with cr_by = "superadmin"
# (a brand-new shared string) instead of "superadmin()".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# regcntr_id, machine_id, device_id triples for the new rows (147-156)
$newRows = @(
    @(10001, 10030, 3000166),
    @(10001, 10030, 3000167),
    @(10001, 10030, 3000168),
    @(10001, 10030, 3000169),
    @(10001, 10030, 3000170),
    @(10001, 10031, 3000171),
    @(10001, 10031, 3000172),
    @(10001, 10031, 3000173),
    @(10001, 10031, 3000174),
    @(10001, 10031, 3000175)
)

$startRow = 147
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = "eng"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin"
    $ws.Cells.Item($r, 7).Value = "now()"
}

# Match the author's final viewport/selection state
$ws.Application.ActiveWindow.ScrollRow = 144
$ws.Range("H149").Select() | Out-Null
